$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1.543815
$ws.Range("A3").Value = 0.371455
$ws.Range("A5").Value = 0.43149
$ws.Range("A7").Value = 0.935825
$ws.Range("A8").Value = 0.37427
$ws.Range("A9").Value = 0.09973499999999999
$ws.Range("A10").Value = 0.097335
$ws.Range("A11").Value = 0.13545
$ws.Range("A12").Value = 0.00755775
$ws.Range("A13").Value = 0.00846015
$ws.Range("A14").Value = 455.645
$ws.Range("A15").Value = 10832.175
